# Scheduled market-data refresh: update currentAveragePrice* / Leve*
# price & profit columns (H:N) for the rows whose figures changed on
# this run. Cells that the refresh could not (re)compute are cleared
# (set to $null) so they serialize as absent, matching the rest of
# the sheet's convention for "no data".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -770

$ws.Range("H33").Value = 140.33333
$ws.Range("I33").Value = 120.35714
$ws.Range("J33").Value = 420
$ws.Range("K33").Value = 120.35714
$ws.Range("L33").Value = 420
$ws.Range("M33").Value = 108.64286
$ws.Range("N33").Value = -878

$ws.Range("H51").Value = 19862.736
$ws.Range("I51").Value = 16166.667
$ws.Range("J51").Value = 20555.75
$ws.Range("K51").Value = 16166.667
$ws.Range("L51").Value = 20555.75
$ws.Range("M51").Value = -15682.667
$ws.Range("N51").Value = -21523.75

$ws.Range("H138").Value = 4528.5713
$ws.Range("I138").Value = 3600
$ws.Range("J138").Value = 6200
$ws.Range("K138").Value = 10800
$ws.Range("L138").Value = 18600
$ws.Range("M138").Value = -5660
$ws.Range("N138").Value = -28880

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 437.125
$ws.Range("I4").Value = 438
$ws.Range("J4").Value = 433.33334
$ws.Range("K4").Value = 438
$ws.Range("L4").Value = 433.33334
$ws.Range("M4").Value = -322
$ws.Range("N4").Value = -665.33334

$ws.Range("H45").Value = 2984
$ws.Range("I45").Value = 2211
$ws.Range("J45").Value = 3499.3333
$ws.Range("K45").Value = 2211
$ws.Range("L45").Value = 3499.3333
$ws.Range("M45").Value = -1834
$ws.Range("N45").Value = -4253.3333

$ws.Range("H97").Value = 357.92307
$ws.Range("I97").Value = 357.92307
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 357.92307
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 138.07693
$ws.Range("N97").Value = $null

$ws.Range("H122").Value = 479664.28
$ws.Range("I122").Value = 669570
$ws.Range("J122").Value = 4900
$ws.Range("K122").Value = 2008710
$ws.Range("L122").Value = 14700
$ws.Range("M122").Value = -2006260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2499.6667
$ws.Range("I86").Value = 2499.6667
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2499.6667
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1376.6667
$ws.Range("N86").Value = $null

$ws.Range("H89").Value = 2499.6667
$ws.Range("I89").Value = 2499.6667
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12498.3335
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -6882.333500000001
$ws.Range("N89").Value = $null

$ws.Range("H94").Value = 837.25
$ws.Range("I94").Value = 424.5
$ws.Range("J94").Value = 1250
$ws.Range("K94").Value = 424.5
$ws.Range("L94").Value = 1250
$ws.Range("M94").Value = 26.5
$ws.Range("N94").Value = -2152

$ws.Range("H99").Value = 1926.8462
$ws.Range("I99").Value = 1575
$ws.Range("J99").Value = 3099.6667
$ws.Range("K99").Value = 1575
$ws.Range("L99").Value = 3099.6667
$ws.Range("M99").Value = -77
$ws.Range("N99").Value = -6095.6667

$ws.Range("H107").Value = 1823.1666
$ws.Range("I107").Value = 1307.8
$ws.Range("J107").Value = 4400
$ws.Range("K107").Value = 1307.8
$ws.Range("L107").Value = 4400
$ws.Range("M107").Value = 612.2
$ws.Range("N107").Value = -8240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 18183138
$ws.Range("I2").Value = 25000802
$ws.Range("J2").Value = 2699.6667
$ws.Range("K2").Value = 25000802
$ws.Range("L2").Value = 2699.6667
$ws.Range("M2").Value = -25000689
$ws.Range("N2").Value = -2925.6667

$ws.Range("H26").Value = 8500
$ws.Range("I26").Value = 2000
$ws.Range("J26").Value = 15000
$ws.Range("K26").Value = 2000
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = -1713
$ws.Range("N26").Value = -15574

$ws.Range("H31").Value = 3969.0344
$ws.Range("I31").Value = 5207.091
$ws.Range("J31").Value = 3212.4443
$ws.Range("K31").Value = 5207.091
$ws.Range("L31").Value = 3212.4443
$ws.Range("M31").Value = -4912.091
$ws.Range("N31").Value = -3802.4443

$ws.Range("H34").Value = 3969.0344
$ws.Range("I34").Value = 5207.091
$ws.Range("J34").Value = 3212.4443
$ws.Range("K34").Value = 5207.091
$ws.Range("L34").Value = 3212.4443
$ws.Range("M34").Value = -5005.091
$ws.Range("N34").Value = -3616.4443

$ws.Range("H132").Value = 2917.8572
$ws.Range("I132").Value = 2477.4546
$ws.Range("J132").Value = 4532.6665
$ws.Range("K132").Value = 7432.3638
$ws.Range("L132").Value = 13597.9995
$ws.Range("M132").Value = -4902.3638
$ws.Range("N132").Value = -18657.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = $null
$ws.Range("N68").Value = $null

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = $null
$ws.Range("N71").Value = $null

$ws.Range("H114").Value = 696
$ws.Range("I114").Value = 751
$ws.Range("J114").Value = 531
$ws.Range("K114").Value = 2253
$ws.Range("L114").Value = 1593
$ws.Range("M114").Value = 1001
$ws.Range("N114").Value = -8101

$ws.Range("H131").Value = 2175.6033
$ws.Range("I131").Value = 1895
$ws.Range("J131").Value = 2184.8032
$ws.Range("K131").Value = 5685
$ws.Range("L131").Value = 6554.409599999999
$ws.Range("M131").Value = -645
$ws.Range("N131").Value = -16634.4096

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = $null

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = $null

$ws.Range("H102").Value = 4079.8
$ws.Range("I102").Value = 1599.6666
$ws.Range("J102").Value = 5142.7144
$ws.Range("K102").Value = 1599.6666
$ws.Range("L102").Value = 5142.7144
$ws.Range("M102").Value = 22.33339999999998

$ws.Range("H107").Value = 1764.64
$ws.Range("I107").Value = 1579.8422
$ws.Range("J107").Value = 2349.8333
$ws.Range("K107").Value = 1579.8422
$ws.Range("L107").Value = 2349.8333
$ws.Range("M107").Value = 340.1578

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 800
$ws.Range("I81").Value = 800
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1600
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -539

$ws.Range("H84").Value = 800
$ws.Range("I84").Value = 800
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 8000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -2696

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = $null

$ws.Range("H132").Value = 5389.364
$ws.Range("I132").Value = 762
$ws.Range("J132").Value = 13487.25
$ws.Range("K132").Value = 2286
$ws.Range("L132").Value = 40461.75
$ws.Range("M132").Value = 244
$ws.Range("N132").Value = -45521.75
